$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Định hướng" column (old column D) entirely -- this shifts
# "Đề tài" and "Điểm" one column to the left (now D, E) and drops the old
# trailing (already empty) column F.
$ws.Range("D1").EntireColumn.Delete()

# Header row (A1:E1) stays: STT, Tên giáo viên, Nhóm, Đề tài, Điểm
$ws.Range("A1").Value = "STT"
$ws.Range("B1").Value = "Tên giáo viên"
$ws.Range("C1").Value = "Nhóm"
$ws.Range("D1").Value = "Đề tài"
$ws.Range("E1").Value = "Điểm"

# Row 2 keeps an (empty) score cell already in place from the column
# shift above -- just update the other columns.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "giaovien1"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "detai1"

# New row 3.
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "giaovien1"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "detai1nhom2"
$ws.Range("E3").ClearFormats()

# New row 4.
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "giaovien2"
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = "nothin"
$ws.Range("E4").ClearFormats()
